$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new text value would otherwise be auto-parsed
# by Excel as a number, so they are stored as literal text (matches the
# source diff, which keeps these as inline/shared strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "76.308.95"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.913.02"
$ws.Range("E3").Value = "  +8.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "196.76"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").Value = "601.97"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "2.910.03"
$ws.Range("E10").Value = "  +8.28%  "
$ws.Range("E11").Value = "  +10.96%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "3.435.12"
$ws.Range("E14").Value = "  +7.77%  "
$ws.Range("D15").Value = "76.203.64"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "2.907.01"
$ws.Range("E18").Value = "  +8.37%  "
$ws.Range("D19").Value = "9.02"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").Value = "385.66"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "72.21"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").Value = "  +7.63%  "
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").Value = "0.0000110"
$ws.Range("E29").Value = "  +15.71%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "515.75"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "7.86"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "165.50"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "20.28"
$ws.Range("E37").Value = "  +5.12%  "
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").Value = "184.22"
$ws.Range("E40").Value = "  +7.78%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "0.349"
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  +10.64%  "
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("D47").Value = "40.26"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "0.588"
$ws.Range("E49").Value = "  +9.04%  "
$ws.Range("D50").Value = "0.689"
$ws.Range("E50").Value = "  +16.26%  "
$ws.Range("D51").Value = "3.79"
$ws.Range("E51").Value = "  +3.44%  "
